# Append new listings and refresh the "取得日時" timestamp on the
# "ランサーズ" sheet (sheet1 / first worksheet in the workbook).
#
# Net effect vs. the original data:
#   - All timestamps in column A move from 2025-12-21 06:35:53
#     to 2025-12-21 12:35:21
#   - Three brand-new listings are merged into the (score-sorted) list:
#       G4 "報酬計算の自動化..."      (score 88)
#       G7 "Manusアプリ..."          (score 38)
#       G9 "グーグルワークスペース..." (score 10)
#   - The whole sheet is rewritten/re-sorted by column G (優先度スコア) desc,
#     growing from 5 data rows (A1:H6) to 8 data rows (A1:H9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-12-21 12:35:21"

# Columns A-H, one row per array entry.
$titles = New-Object 'object[]' 8
$titles[0] = "【急募】ECサイトの自動購入Bot作成をお願いします。"
$titles[1] = "【受注メールを元にECサイト自動仕入ツール】"
$titles[2] = "【報酬計算の自動化】GASで自動計算させるプログラミング"
$titles[3] = "【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー"
$titles[4] = "初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)"
$titles[5] = "【急募】ManusアプリのGoogleStore登録代行を依頼したい"
$titles[6] = "【急募】Notionでの社内向けダッシュボード作成依頼"
$titles[7] = "【急募】グーグルワークスペースの設定をサポートしてくれる方"

$categories = New-Object 'object[]' 8
for ($i = 0; $i -lt 8; $i++) { $categories[$i] = "システム開発" }

$prices = New-Object 'object[]' 8
$prices[0] = "100,000 円 ~ 200,000 円 / 固定"
$prices[1] = "50,000 円 ~ 100,000 円 / 固定"
$prices[2] = "50,000 円 ~ 100,000 円 / 固定"
$prices[3] = "200,000 円 ~ 300,000 円 / 固定"
$prices[4] = "300,000 円 ~ 500,000 円 / 固定"
$prices[5] = "50,000 円 ~ 100,000 円 / 固定"
$prices[6] = "50,000 円 ~ 100,000 円 / 固定"
$prices[7] = "1,000 ~ 5,000 円 / 固定"

$deadlines = New-Object 'object[]' 8
for ($i = 0; $i -lt 8; $i++) { $deadlines[$i] = "期限情報なし" }

$urls = New-Object 'object[]' 8
$urls[0] = "https://www.lancers.jp/work/detail/5458190"
$urls[1] = "https://www.lancers.jp/work/detail/5458166"
$urls[2] = "https://www.lancers.jp/work/detail/5458299"
$urls[3] = "https://www.lancers.jp/work/detail/5431107"
$urls[4] = "https://www.lancers.jp/work/detail/5425629"
$urls[5] = "https://www.lancers.jp/work/detail/5458330"
$urls[6] = "https://www.lancers.jp/work/detail/5458234"
$urls[7] = "https://www.lancers.jp/work/detail/5458288"

$scores = New-Object 'object[]' 8
$scores[0] = 143
$scores[1] = 98
$scores[2] = 88
$scores[3] = 68
$scores[4] = 45
$scores[5] = 38
$scores[6] = 18
$scores[7] = 10

$skills = New-Object 'object[]' 8
$skills[0] = "★bot ◇サイト"
$skills[1] = "◆ツール ◇サイト"
$skills[2] = "◆自動化"
$skills[3] = "◆開発"
$skills[4] = "◇サイト"
$skills[5] = "◇アプリ"
$skills[6] = ""
$skills[7] = ""

$rowCount = 8
$colCount = 8
$data = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $timestamp
    $data[$i,1] = $titles[$i]
    $data[$i,2] = $categories[$i]
    $data[$i,3] = $prices[$i]
    $data[$i,4] = $deadlines[$i]
    $data[$i,5] = $urls[$i]
    $data[$i,6] = $scores[$i]
    $data[$i,7] = $skills[$i]
}

# Write the whole block (rows 2-9, columns A-H) in one shot.
$ws.Range("A2:H9").Value = $data

# Clear H8 / H9 (skill summary is blank for those two listings, matching
# the source rows which never had an H cell at all).
$ws.Range("H8").Value = ""
$ws.Range("H9").Value = ""

# Drop the hyperlinks that were already attached to the old F2:F6 range so
# they don't linger as stale/duplicate links once the rows are rewritten.
$ws.Range("F2:F9").Hyperlinks.Delete()

# Re-create the hyperlinks on column F for every data row, and make sure
# they keep using the workbook's "Hyperlink" cell style (as the original
# rows F2:F6 already did).
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Range("F" + $r)
    $ws.Hyperlinks.Add($cell, $urls[$r - 2])
    $cell.Style = "Hyperlink"
}

$ws.Range("A1").Select()
